$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.953.37'
$ws.Range("E2").Value = '  -1.48%  '
$ws.Range("D3").Value = '1.892.15'
$ws.Range("E3").Value = '  -2.34%  '
$ws.Range("E4").Value = '  -0.33%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7347'
$ws.Range("E5").Value = '  -1.65%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '242.72'
$ws.Range("E6").Value = '  -1.01%  '
$ws.Range("E7").Value = '  -0.31%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3090'
$ws.Range("E8").Value = '  -2.53%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '26.33'
$ws.Range("E9").Value = '  -4.48%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06895'
$ws.Range("E10").Value = '  -1.23%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7700'
$ws.Range("E11").Value = '  -1.31%  '
$ws.Range("E12").Value = '  -0.55%  '
$ws.Range("D13").Value = '1.901.99'
$ws.Range("E13").Value = '  -1.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.212'
$ws.Range("E14").Value = '  -2.84%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.39'
$ws.Range("E15").Value = '  -3.25%  '
$ws.Range("D16").Value = '29.967.45'
$ws.Range("E16").Value = '  -1.39%  '
$ws.Range("E17").Value = '  -2.22%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.788'
$ws.Range("E18").Value = '  +0.85%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '239.26'
$ws.Range("E19").Value = '  -5.39%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007757'
$ws.Range("E20").Value = '  -2.11%  '
$ws.Range("E21").Value = '  -0.22%  '
$ws.Range("D22").Value = '2.136.52'
$ws.Range("E22").Value = '  -2.61%  '
$ws.Range("E23").Value = '  -0.34%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.951'
$ws.Range("E24").Value = '  +4.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.283'
$ws.Range("E25").Value = '  -2.36%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.16'
$ws.Range("E26").Value = '  -0.45%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.78'
$ws.Range("E27").Value = '  -1.08%  '
$ws.Range("E28").Value = '  -4.56%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.021'
$ws.Range("E29").Value = '  -9.85%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.354'
$ws.Range("E30").Value = '  -0.74%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.535'
$ws.Range("E31").Value = '  +1.63%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.294'
$ws.Range("E32").Value = '  -1.48%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.055'
$ws.Range("E33").Value = '  -1.30%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05089'
$ws.Range("E34").Value = '  -1.26%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.273'
$ws.Range("E35").Value = '  +0.05%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7334'
$ws.Range("E36").Value = '  -1.83%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.718'
$ws.Range("E37").Value = '  -2.47%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01921'
$ws.Range("E38").Value = '  -1.46%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.784'
$ws.Range("E39").Value = '  -0.55%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.296'
$ws.Range("E40").Value = '  -1.87%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '74.13'
$ws.Range("E41").Value = '  -4.82%  '
$ws.Range("E42").Value = '  -0.37%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.929'
$ws.Range("E43").Value = '  -1.69%  '
$ws.Range("E44").Value = '  -0.26%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8353'
$ws.Range("E45").Value = '  +0.32%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.634'
$ws.Range("E46").Value = '  +2.50%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '101.00'
$ws.Range("E47").Value = '  -0.37%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.776'
$ws.Range("E48").Value = '  -0.03%  '
$ws.Range("D49").Value = '2.038.73'
$ws.Range("E49").Value = '  -2.34%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.36'
$ws.Range("E50").Value = '  -2.38%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '932.81'
$ws.Range("E51").Value = '  -5.34%  '
